$wb = $excel.ActiveWorkbook

# Property sheet: add new last column T1 = property__topic
$wsProperty = $wb.Worksheets.Item("Property")
$wsProperty.Range("T1").Value = "property__topic"

# SMWType sheet: add new last column G1 = sMWType__usedByProperties
$wsSMWType = $wb.Worksheets.Item("SMWType")
$wsSMWType.Range("G1").Value = "sMWType__usedByProperties"

# Topic sheet: add new last column K1 = topic__context
$wsTopic = $wb.Worksheets.Item("Topic")
$wsTopic.Range("K1").Value = "topic__context"
